# Apply updated loading_percent values for Case_2_250 (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.86867074118826
$ws.Range("C2").Value = 12.44650364804212
$ws.Range("E2").Value = 16.88878898332405
$ws.Range("F2").Value = 38.1279623031242
$ws.Range("G2").Value = 32.33038485469088
$ws.Range("H2").Value = 15.35259844640018
$ws.Range("J2").Value = 8.058749281693816
$ws.Range("K2").Value = 8.389170202324598
$ws.Range("L2").Value = 12.35103032800004
$ws.Range("O2").Value = 23.81053906310736

$ws.Range("B3").Value = 12.5863359163747
$ws.Range("C3").Value = 12.48932229835657
$ws.Range("E3").Value = 16.88757236665566
$ws.Range("F3").Value = 38.18137056632427
$ws.Range("G3").Value = 32.48996141659094
$ws.Range("H3").Value = 15.4082909874739
$ws.Range("J3").Value = 8.051377812484569
$ws.Range("K3").Value = 8.175709224834742
$ws.Range("L3").Value = 12.32093068236677
$ws.Range("O3").Value = 23.91279849931067

$ws.Range("B4").Value = 12.41127141871184
$ws.Range("C4").Value = 12.51698053927911
$ws.Range("E4").Value = 16.88950195653977
$ws.Range("F4").Value = 38.22280680706417
$ws.Range("G4").Value = 32.59680555999698
$ws.Range("H4").Value = 15.44466551970108
$ws.Range("J4").Value = 8.046892837478063
$ws.Range("K4").Value = 8.042327645248399
$ws.Range("L4").Value = 12.30405346897593
$ws.Range("O4").Value = 23.98005412457545

$ws.Range("B5").Value = 12.33959792316548
$ws.Range("C5").Value = 12.52859632856396
$ws.Range("E5").Value = 16.89096311289856
$ws.Range("F5").Value = 38.24186369294289
$ws.Range("G5").Value = 32.64256725395901
$ws.Range("H5").Value = 15.4600369821835
$ws.Range("J5").Value = 8.045075809101967
$ws.Range("K5").Value = 7.987464238387264
$ws.Range("L5").Value = 12.29758383371534
$ws.Range("O5").Value = 24.0085844922321

$ws.Range("B6").Value = 12.32767935704076
$ws.Range("C6").Value = 12.53054597747919
$ws.Range("E6").Value = 16.89124654196696
$ws.Range("F6").Value = 38.24515913694894
$ws.Range("G6").Value = 32.6502999376061
$ws.Range("H6").Value = 15.46262255037555
$ws.Range("J6").Value = 8.044774743703227
$ws.Range("K6").Value = 7.978325660928165
$ws.Range("L6").Value = 12.29653432448002
$ws.Range("O6").Value = 24.01338976447531

$ws.Range("B7").Value = 12.4103060240796
$ws.Range("C7").Value = 12.51713579635585
$ws.Range("E7").Value = 16.88951892766973
$ws.Range("F7").Value = 38.22305502685403
$ws.Range("G7").Value = 32.59741373130808
$ws.Range("H7").Value = 15.44487060278306
$ws.Range("J7").Value = 8.046868288955494
$ws.Range("K7").Value = 8.041589701022952
$ws.Range("L7").Value = 12.30396455951206
$ws.Range("O7").Value = 23.98043434791408

$ws.Range("B8").Value = 12.77173340181784
$ws.Range("C8").Value = 12.46098430457061
$ws.Range("E8").Value = 16.88781510771865
$ws.Range("F8").Value = 38.14458245705418
$ws.Range("G8").Value = 32.38356307188556
$ws.Range("H8").Value = 15.37134932316374
$ws.Range("J8").Value = 8.056199235280495
$ws.Range("K8").Value = 8.316091970830719
$ws.Range("L8").Value = 12.34032188998235
$ws.Range("O8").Value = 23.8448705935278

$ws.Range("B9").Value = 13.46287598157538
$ws.Range("C9").Value = 12.36167615916206
$ws.Range("E9").Value = 16.90562418452457
$ws.Range("F9").Value = 38.05934770245335
$ws.Range("G9").Value = 32.03483878546146
$ws.Range("H9").Value = 15.24443729313907
$ws.Range("J9").Value = 8.07481139628624
$ws.Range("K9").Value = 8.833025890468223
$ws.Range("L9").Value = 12.42413236361643
$ws.Range("O9").Value = 23.61449654093698

$ws.Range("B10").Value = 13.95481054343059
$ws.Range("C10").Value = 12.2952379599454
$ws.Range("E10").Value = 16.93147162803886
$ws.Range("F10").Value = 38.03862849480254
$ws.Range("G10").Value = 31.82211040973855
$ws.Range("H10").Value = 15.16168124303092
$ws.Range("J10").Value = 8.088661786725723
$ws.Range("K10").Value = 9.196105574029339
$ws.Range("L10").Value = 12.49303361370532
$ws.Range("O10").Value = 23.46688245409512

$ws.Range("B11").Value = 14.1741802180173
$ws.Range("C11").Value = 12.26641656379449
$ws.Range("E11").Value = 16.94596831972532
$ws.Range("F11").Value = 38.03829618045638
$ws.Range("G11").Value = 31.73486406552745
$ws.Range("H11").Value = 15.12630261830107
$ws.Range("J11").Value = 8.094997568503661
$ws.Range("K11").Value = 9.356963864501832
$ws.Range("L11").Value = 12.52589895019781
$ws.Range("O11").Value = 23.40443344440409

$ws.Range("B12").Value = 14.25654026537071
$ws.Range("C12").Value = 12.25570320581388
$ws.Range("E12").Value = 16.95184831217986
$ws.Range("F12").Value = 38.03947589765795
$ws.Range("G12").Value = 31.7032030366932
$ws.Range("H12").Value = 15.11323119970688
$ws.Range("J12").Value = 8.097401503830485
$ws.Range("K12").Value = 9.417206105254209
$ws.Range("L12").Value = 12.53855662417105
$ws.Range("O12").Value = 23.38146229310198

$ws.Range("B13").Value = 14.23883530363017
$ws.Range("C13").Value = 12.25800160778773
$ws.Range("E13").Value = 16.95056464401152
$ws.Range("F13").Value = 38.03916380072457
$ws.Range("G13").Value = 31.70996043150179
$ws.Range("H13").Value = 15.11603188369543
$ws.Range("J13").Value = 8.096883569318116
$ws.Range("K13").Value = 9.404262492842504
$ws.Range("L13").Value = 12.53582123315338
$ws.Range("O13").Value = 23.38637942037302

$ws.Range("B14").Value = 14.18097068269451
$ws.Range("C14").Value = 12.26553115247998
$ws.Range("E14").Value = 16.94644426627345
$ws.Range("F14").Value = 38.03836708579305
$ws.Range("G14").Value = 31.73223164797231
$ws.Range("H14").Value = 15.12522069937786
$ws.Range("J14").Value = 8.095195249264682
$ws.Range("K14").Value = 9.36193376404815
$ws.Range("L14").Value = 12.52693608856593
$ws.Range("O14").Value = 23.40253002204666

$ws.Range("B15").Value = 14.14543224860405
$ws.Range("C15").Value = 12.27016932442822
$ws.Range("E15").Value = 16.94397115336856
$ws.Range("F15").Value = 38.03804902252325
$ws.Range("G15").Value = 31.74605299586093
$ws.Range("H15").Value = 15.13089151930487
$ws.Range("J15").Value = 8.094161707825814
$ws.Range("K15").Value = 9.335917313945199
$ws.Range("L15").Value = 12.52152112795971
$ws.Range("O15").Value = 23.41251091531948

$ws.Range("B16").Value = 13.94037879967784
$ws.Range("C16").Value = 12.29714963580593
$ws.Range("E16").Value = 16.93057905848751
$ws.Range("F16").Value = 38.03883303527779
$ws.Range("G16").Value = 31.8280044807054
$ws.Range("H16").Value = 15.16403891577863
$ws.Range("J16").Value = 8.088248406814817
$ws.Range("K16").Value = 9.185501880799901
$ws.Range("L16").Value = 12.49091585524946
$ws.Range("O16").Value = 23.47105832204259

$ws.Range("B17").Value = 13.81339902209668
$ws.Range("C17").Value = 12.31405955625776
$ws.Range("E17").Value = 16.92306252948277
$ws.Range("F17").Value = 38.04164178542729
$ws.Range("G17").Value = 31.88072405350282
$ws.Range("H17").Value = 15.18495426825084
$ws.Range("J17").Value = 8.08462955216843
$ws.Range("K17").Value = 9.092085104496677
$ws.Range("L17").Value = 12.47252570461737
$ws.Range("O17").Value = 23.50818001179893

$ws.Range("B18").Value = 13.73995235083641
$ws.Range("C18").Value = 12.32391767942441
$ws.Range("E18").Value = 16.91899720795212
$ws.Range("F18").Value = 38.04411335194758
$ws.Range("G18").Value = 31.91194303415419
$ws.Range("H18").Value = 15.19719767608669
$ws.Range("J18").Value = 8.082551414007474
$ws.Range("K18").Value = 9.037951510137813
$ws.Range("L18").Value = 12.46209182504796
$ws.Range("O18").Value = 23.52997388713205

$ws.Range("B19").Value = 13.71501633337336
$ws.Range("C19").Value = 12.32727816591671
$ws.Range("E19").Value = 16.91766516819096
$ws.Range("F19").Value = 38.04509725531522
$ws.Range("G19").Value = 31.92266691312947
$ws.Range("H19").Value = 15.20137975540639
$ws.Range("J19").Value = 8.081848374717961
$ws.Range("K19").Value = 9.019555286866265
$ws.Range("L19").Value = 12.45858395993149
$ws.Range("O19").Value = 23.5374288857116

$ws.Range("B20").Value = 13.82695936641215
$ws.Range("C20").Value = 12.31224581289431
$ws.Range("E20").Value = 16.9238359998026
$ws.Range("F20").Value = 38.04125420310995
$ws.Range("G20").Value = 31.87501916894058
$ws.Range("H20").Value = 15.18270570468609
$ws.Range("J20").Value = 8.085014443126322
$ws.Range("K20").Value = 9.102071550750647
$ws.Range("L20").Value = 12.47446854709997
$ws.Range("O20").Value = 23.50418254519283

$ws.Range("B21").Value = 14.19798678399401
$ws.Range("C21").Value = 12.26331410327229
$ws.Range("E21").Value = 16.94764395390642
$ws.Range("F21").Value = 38.03856568816433
$ws.Range("G21").Value = 31.72565261313068
$ws.Range("H21").Value = 15.12251288493853
$ws.Range("J21").Value = 8.095691024676913
$ws.Range("K21").Value = 9.374385346571989
$ws.Range("L21").Value = 12.52954016301529
$ws.Range("O21").Value = 23.39776781987887

$ws.Range("B22").Value = 14.43630205951268
$ws.Range("C22").Value = 12.2325037130595
$ws.Range("E22").Value = 16.96547766478685
$ws.Range("F22").Value = 38.04441681934171
$ws.Range("G22").Value = 31.63606415174814
$ws.Range("H22").Value = 15.08507166177206
$ws.Range("J22").Value = 8.102696090756293
$ws.Range("K22").Value = 9.548424462115747
$ws.Range("L22").Value = 12.56676679197999
$ws.Range("O22").Value = 23.33216576283636

$ws.Range("B23").Value = 14.30951413888896
$ws.Range("C23").Value = 12.24884109436898
$ws.Range("E23").Value = 16.95575261449529
$ws.Range("F23").Value = 38.04059869732106
$ws.Range("G23").Value = 31.68314186055763
$ws.Range("H23").Value = 15.10488117260293
$ws.Range("J23").Value = 8.098954972340401
$ws.Range("K23").Value = 9.455912371284093
$ws.Range("L23").Value = 12.54678754287977
$ws.Range("O23").Value = 23.36681744575878

$ws.Range("B24").Value = 13.8208301128421
$ws.Range("C24").Value = 12.31306538137474
$ws.Range("E24").Value = 16.92348551611338
$ws.Range("F24").Value = 38.0414267601456
$ws.Range("G24").Value = 31.87759551386743
$ws.Range("H24").Value = 15.18372159840195
$ws.Range("J24").Value = 8.084840426453084
$ws.Range("K24").Value = 9.097558004573203
$ws.Range("L24").Value = 12.47358975498074
$ws.Range("O24").Value = 23.50598839157406

$ws.Range("B25").Value = 13.27833569586052
$ws.Range("C25").Value = 12.3873915062945
$ws.Range("E25").Value = 16.89855446011488
$ws.Range("F25").Value = 38.0750479285633
$ws.Range("G25").Value = 32.1215724575635
$ws.Range("H25").Value = 15.27692621000383
$ws.Range("J25").Value = 8.069744256048491
$ws.Range("K25").Value = 8.695867449574679
$ws.Range("L25").Value = 12.40014984699404
$ws.Range("O25").Value = 23.67301993843298
